$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 21 ("Pour toi le beur" / large bibliography text).
# This shifts rows 22-24 up by one, so:
#   old row22 (Pousse Les Bass / ".")      -> new row21
#   old row23 (Rappeur / animalerie text)  -> new row22
#   old row24 (Y'a pas de problème / lyrics)-> new row23
# and the sheet shrinks from 24 to 23 rows.
$ws.Rows.Item(21).Delete()
